# Updates cryptos list data (prices / 1h volume %) and fixes the ordering
# of a few rows (Bittensor/Cosmos, Stacks/InjectiveProtocol/PEPE,
# THORChain/Monero), matching the GitHub Actions data refresh commit.
#
# Price cells (column D) hold exact text (e.g. "6.01", "1.00", "64.826.48")
# rather than numbers, so a leading "'" (quote-prefix) is used wherever the
# new value would otherwise be auto-parsed as a number by Excel - this keeps
# the cell as text and preserves formatting such as trailing zeros without
# introducing floating point rounding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.826.48"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").Value = "3.433.56"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'570.92"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'174.39"
$ws.Range("E6").Value = "  -7.74%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("D11").Value = "'54.94"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "'9.12"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").Value = "3.983.96"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "3.442.36"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "64.815.53"
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("D19").Value = "'11.82"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "'0.989"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'407.03"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "'4.32"
$ws.Range("E23").Value = "  +4.66%  "
$ws.Range("D24").Value = "'83.43"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "'13.27"
$ws.Range("E25").Value = "  +8.46%  "
$ws.Range("D26").Value = "'10.81"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").Value = "'6.01"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "'8.95"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "'29.85"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").Value = "'6.64"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'590.40"
$ws.Range("E32").Value = "  -8.44%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.51"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "'59.48"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'36.10"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0767"
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.52"
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("D41").Value = "'0.376"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").Value = "3.180.97"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "'2.91"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "'2.51"
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "'2.64"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'137.99"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.43"
$ws.Range("E51").Value = "  -2.16%  "
